# FWHM data run for sg_rr_100_030 2023-12-08 16-08-32
# Append a new data row (row 62) to the FWHM results table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

$ws.Cells.Item($row, 1).Value  = "sg_rr_100_030 2023-12-08 16-08-32.csv"   # Data CSV Filename
$ws.Cells.Item($row, 2).Value  = 0.01                                     # Wavelength step size/nm
$ws.Cells.Item($row, 3).Value  = 1000                                     # Start array index
$ws.Cells.Item($row, 4).Value  = 5001                                     # End array index
$ws.Cells.Item($row, 5).Value  = 1530                                     # Start wavelength/nm
$ws.Cells.Item($row, 6).Value  = 1570                                     # End wavelength/nm
$ws.Cells.Item($row, 7).Value  = 0.5                                      # prominence/dBm
$ws.Cells.Item($row, 8).Value  = "(approx_fsr/2)/wavelength step size"    # distance
$ws.Cells.Item($row, 9).Value  = 1                                        # approx_fsr/nm
$ws.Cells.Item($row, 10).Value = 0.98128205128205004                      # fsr_mean/nm
$ws.Cells.Item($row, 11).Value = 0.00316397329552258                      # fsr_std error/nm
$ws.Cells.Item($row, 12).Value = "yes"                                    # double count check passed?
$ws.Cells.Item($row, 13).Value = 0.0780266292280477                       # mean FWHM/nm
$ws.Cells.Item($row, 14).Value = 0.00196369374415618                      # FWHM error/nm

# Scroll the view up a bit and leave the selection on the new FWHM-error cell.
$ws.Application.ActiveWindow.ScrollRow = 41
$ws.Range("N62").Select()
